$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet previously had an AutoFilter applied on column A (value "Cortes"),
# which hid every row whose A-column didn't match. Clear the filter so every
# row becomes visible again (also drops the <filterColumn>/<filters> entry
# and clears sheetPr's filterMode where the host supports it).
$ws.ShowAllData()

# Scroll the view down so row 10 becomes the top visible row, then put the
# selection on A21 (single cell), matching the saved view state.
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 10
$ws.Range("A21").Select()
